# Update the Avanade "Analyst" entry: add "(Contractor)" to the title and
# give the role an end date (Present -> October 2018).

$d = $word.ActiveDocument

# 1) "Analyst, Machine Learning and Azure Cloud Enablement"
#    -> "Analyst (Contractor), Machine Learning and Azure Cloud Enablement"
$d.Content.Find.Execute(
    "Analyst, Machine Learning and Azure Cloud Enablement",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Analyst (Contractor), Machine Learning and Azure Cloud Enablement",
    2
)

# 2) "January 2018 – Present" -> "January 2018 – October 2018"
$d.Content.Find.Execute(
    "January 2018 " + [char]8211 + " Present",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "January 2018 " + [char]8211 + " October 2018",
    2
)
